$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Role"
$ws.Range("C4").Value = "UserRole"
$ws.Range("E4").Value = "Movie"
$ws.Range("F4").Value = "Project"
$ws.Range("G4").Value = "Like"
$ws.Range("A4").Value = "User"
$ws.Range("I4").Value = "Subscriber"
$ws.Range("J4").Value = "WantedSubscriber"
$ws.Range("D4").Value = "Picture"
$ws.Range("H4").Value = "Category"

$ws.Range("H5").Select()
